# Update "想去人数" (interested-people count) values in column F
# on the "展览" and "全部类型" worksheets, per the commit's regenerated data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - all 14 rows updated
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 1105
$wsExpo.Range("F5").Value  = 88
$wsExpo.Range("F8").Value  = 11342
$wsExpo.Range("F9").Value  = 4310
$wsExpo.Range("F11").Value = 27
$wsExpo.Range("F12").Value = 16
$wsExpo.Range("F13").Value = 2519
$wsExpo.Range("F14").Value = 1075
$wsExpo.Range("F15").Value = 115
$wsExpo.Range("F16").Value = 21
$wsExpo.Range("F17").Value = 173
$wsExpo.Range("F19").Value = 11265
$wsExpo.Range("F20").Value = 11130
$wsExpo.Range("F25").Value = 38

# Sheet "全部类型" (All Types) - same rows, except F11 was already 27
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 1105
$wsAll.Range("F5").Value  = 88
$wsAll.Range("F8").Value  = 11342
$wsAll.Range("F9").Value  = 4310
$wsAll.Range("F12").Value = 16
$wsAll.Range("F13").Value = 2519
$wsAll.Range("F14").Value = 1075
$wsAll.Range("F15").Value = 115
$wsAll.Range("F16").Value = 21
$wsAll.Range("F17").Value = 173
$wsAll.Range("F19").Value = 11265
$wsAll.Range("F20").Value = 11130
$wsAll.Range("F25").Value = 38
